$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.349.09"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "'1.868.13"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'243.65"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.4702"
$ws.Range("E7").Value = "  -1.18%  "
$ws.Range("E8").Value = "  -1.16%  "
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("D11").Value = "'0.07766"
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").Value = "'1.867.34"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").Value = "'95.82"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").Value = "'0.7203"
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").Value = "'5.124"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").Value = "'278.75"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "'30.340.24"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").Value = "'12.96"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").Value = "'0.000007527"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "'1.0000"
$ws.Range("D21").Value = "'2.115.59"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'5.219"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "'6.221"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").Value = "'163.41"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("E26").Value = "  -1.48%  "
$ws.Range("D27").Value = "'18.64"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("D29").Value = "'1.320"
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("D30").Value = "'0.09593"
$ws.Range("E30").Value = "  -2.73%  "
$ws.Range("D31").Value = "'1.467"
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("D32").Value = "'4.201"
$ws.Range("D33").Value = "'4.096"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").Value = "'0.04806"
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("D35").Value = "'1.117"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "'0.6884"
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("D37").Value = "'2.713"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "'0.01874"
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("D40").Value = "'6.206"
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("D41").Value = "'74.12"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").Value = "'1.937"
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("D43").Value = "'0.4213"
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "'0.8244"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("D46").Value = "'100.77"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").Value = "'9.561"
$ws.Range("E47").Value = "  +2.25%  "
$ws.Range("D48").Value = "'35.23"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").Value = "'6.944"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").Value = "'898.58"
$ws.Range("E50").Value = "  -1.82%  "
